$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Text tweaks on "LinaJourney" (sheet1)
# ---------------------------------------------------------------------------
$journey = $wb.Worksheets.Item("LinaJourney")
$journey.Range("B8").Value  = "jesus ✨"
$journey.Range("B10").Value = "gimme a minute 🎶"
$journey.Range("B11").Value = "I'm done holding back"
$journey.Range("B20").Value = "well? ✨"

# ---------------------------------------------------------------------------
# 2) "cumcontrol" becomes "cumcontrol1" with refreshed copy / sync / edge
#    lines (moving toward explicit "send PPV" / dom-control language).
# ---------------------------------------------------------------------------
$cc = $wb.Worksheets.Item("cumcontrol")
$cc.Name = "cumcontrol1"

$cc.Range("B2").Value = "trust me you want to edge just a little longer for this one 💕"

$cc.Range("B3").Value = "you're not done until I say you are... open this"
$cc.Range("C3").Value = "DELAY. Send PPV."

$cc.Range("B4").Value = "I'm right there too, let's finish this... but you need to see this first"
$cc.Range("C4").Value = "SYNC variant. Send PPV."

$cc.Range("B5").Value = "now... right now, with me love. open this ✨"
$cc.Range("C5").Value = "SYNC. Send PPV."

$cc.Range("B6").Value = "not a chance... you're going to wait until I say so"

$cc.Range("B7").Value = "I didn't say you could cum yet love ✨"
$cc.Range("C7").Value = "CONTROL."

# ---------------------------------------------------------------------------
# 3) Insert a brand new "cumcontrol2" sheet right after "cumcontrol1" — a
#    duplicate of its layout/styles, re-purposed with delay/sync/edge
#    variant copy (this is what used to live under "dickpic").
#    Doing this as a copy of cumcontrol1 pushes "dickpic" and "boosters"
#    down by one slot, unchanged.
# ---------------------------------------------------------------------------
$cc.Copy($null, $cc)
$cc2 = $wb.Worksheets.Item("cumcontrol1 (2)")
$cc2.Name = "cumcontrol2"

$cc2.Range("A2").Value = "delay2"
$cc2.Range("B2").Value = "edge for me... just a little more... this last one is everything ✨"
$cc2.Range("C2").Value = "DELAY variant."

$cc2.Range("A3").Value = "delay1"
$cc2.Range("B3").Value = "hold it... what I'm about to send is the best one and you'll want to last for it"
$cc2.Range("C3").Value = "DELAY. Send PPV."

$cc2.Range("A4").Value = "sync2"
$cc2.Range("B4").Value = "I want us to finish at the same time... this one will push you over"
$cc2.Range("C4").Value = "SYNC variant."

$cc2.Range("A5").Value = "sync1"
$cc2.Range("B5").Value = "okay you earned it love... let's go together, open this 💕"
$cc2.Range("C5").Value = "SYNC. Send PPV."

$cc2.Range("A6").Value = "edge2"
$cc2.Range("B6").Value = "if you finish without my permission I'll be annoyed"
$cc2.Range("C6").Value = "EDGE variant."

$cc2.Range("A7").Value = "edge1"
$cc2.Range("B7").Value = "slow down love, I'm in control here ✨"
$cc2.Range("C7").Value = "CONTROL."

# Restore original active sheet/selection (copying a sheet activates it).
$journey.Select() | Out-Null
